$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 5946
$ws.Range("J3").Value = 8077
$ws.Range("K3").Value = 6124
$ws.Range("H4").Value = 1739
$ws.Range("K4").Value = 1280
$ws.Range("K5").Value = 435
$ws.Range("J6").Value = 11052
$ws.Range("K6").Value = 6741
$ws.Range("H7").Value = 26052
$ws.Range("J7").Value = 29299
$ws.Range("K7").Value = 20526

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K3").Value = 70
$ws.Range("K7").Value = 262

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K6").Value = 455
$ws.Range("K7").Value = 1358

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 157
$ws.Range("K7").Value = 454

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 327
$ws.Range("K5").Value = 23
$ws.Range("K7").Value = 887

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 230
$ws.Range("K6").Value = 202
$ws.Range("K7").Value = 695

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K3").Value = 117
$ws.Range("K7").Value = 479

$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range("K3").Value = 139
$ws.Range("K7").Value = 337

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K6").Value = 28
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K5").Value = 52
$ws.Range("K7").Value = 596
$ws.Range("K8").Value = 1358
$ws.Range("K15").Value = 211
$ws.Range("K18").Value = 134
$ws.Range("K19").Value = 593
$ws.Range("K20").Value = 486
$ws.Range("K22").Value = 60
$ws.Range("J27").Value = 180
$ws.Range("K29").Value = 1120
$ws.Range("K30").Value = 77
$ws.Range("K33").Value = 887
$ws.Range("K36").Value = 262
$ws.Range("K37").Value = 695
$ws.Range("J42").Value = 1238
$ws.Range("K47").Value = 142
$ws.Range("K48").Value = 257
$ws.Range("K49").Value = 112
$ws.Range("K51").Value = 264
$ws.Range("K53").Value = 262
$ws.Range("K54").Value = 399
$ws.Range("K60").Value = 122
$ws.Range("H61").Value = 34
$ws.Range("K63").Value = 67
$ws.Range("K65").Value = 479
$ws.Range("K67").Value = 802
$ws.Range("K68").Value = 56
$ws.Range("K72").Value = 99
$ws.Range("K76").Value = 278
$ws.Range("K78").Value = 234
$ws.Range("K79").Value = 508
$ws.Range("K81").Value = 16
$ws.Range("K83").Value = 454
$ws.Range("K85").Value = 954
$ws.Range("K88").Value = 221
$ws.Range("K91").Value = 233
$ws.Range("K92").Value = 78
$ws.Range("K93").Value = 76
$ws.Range("K94").Value = 274
$ws.Range("K99").Value = 337
$ws.Range("H101").Value = 26052
$ws.Range("J101").Value = 29299
$ws.Range("K101").Value = 20526

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K3").Value = 285
$ws.Range("K7").Value = 802

$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 112

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K6").Value = 216
$ws.Range("K7").Value = 399

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 319
$ws.Range("K7").Value = 1120

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K3").Value = 62
$ws.Range("K7").Value = 257

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 182
$ws.Range("K4").Value = 27
$ws.Range("K6").Value = 190
$ws.Range("K7").Value = 593

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K3").Value = 53
$ws.Range("K6").Value = 145
$ws.Range("K7").Value = 278

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("J3").Value = 250
$ws.Range("J7").Value = 1238

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 234

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 58
$ws.Range("K4").Value = 10
$ws.Range("K7").Value = 233

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K3").Value = 166
$ws.Range("K6").Value = 124
$ws.Range("K7").Value = 508

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K6").Value = 139
$ws.Range("K7").Value = 486

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K2").Value = 36
$ws.Range("K7").Value = 134

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("K3").Value = 78
$ws.Range("K7").Value = 262

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K6").Value = 29
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K6").Value = 160
$ws.Range("K7").Value = 596

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K2").Value = 73
$ws.Range("K7").Value = 274

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K2").Value = 42
$ws.Range("K7").Value = 142

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Range("K6").Value = 66
$ws.Range("K7").Value = 211

$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range("K3").Value = 15
$ws.Range("K7").Value = 78

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 56
$ws.Range("K7").Value = 221

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 52

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("J6").Value = 64
$ws.Range("J7").Value = 180

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K5").Value = 6
$ws.Range("K7").Value = 264

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("K3").Value = 12
$ws.Range("K7").Value = 56

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 122

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 310
$ws.Range("K3").Value = 329
$ws.Range("K7").Value = 954

$ws = $wb.Worksheets.Item('Clearing')
$ws.Range("K6").Value = 12
$ws.Range("K7").Value = 60

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K2").Value = 17
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item('Mount Greenwood')
$ws.Range("H4").Value = 1
$ws.Range("H7").Value = 34

$ws = $wb.Worksheets.Item('Sauganash,Forest Glen')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 16
